$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 12 and 13 (columns F..V) ---
$ws.Range("F12").Value = "Oran"
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = "Saoura"
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 2.54
$ws.Range("K12").Value = "22/09/2023 08:13"
$ws.Range("L12").Value = 2.34
$ws.Range("M12").Value = "23/09/2023 18:03"
$ws.Range("N12").Value = 2.65
$ws.Range("O12").Value = "22/09/2023 08:13"
$ws.Range("P12").Value = 2.71
$ws.Range("Q12").Value = "23/09/2023 18:03"
$ws.Range("R12").Value = 3.11
$ws.Range("S12").Value = "22/09/2023 08:13"
$ws.Range("T12").Value = 3.84
$ws.Range("U12").Value = "23/09/2023 18:03"
$ws.Range("V12").Value = "https://www.betexplorer.com/football/algeria/ligue-1/oran-saoura/YyyeO358/"

$ws.Range("F13").Value = "Biskra"
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = "US Souf"
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1.71
$ws.Range("K13").Value = "22/09/2023 15:13"
$ws.Range("L13").Value = 1.65
$ws.Range("M13").Value = "23/09/2023 19:13"
$ws.Range("N13").Value = 3.22
$ws.Range("O13").Value = "22/09/2023 15:13"
$ws.Range("P13").Value = 3.47
$ws.Range("Q13").Value = "23/09/2023 19:13"
$ws.Range("R13").Value = 4.87
$ws.Range("S13").Value = "22/09/2023 15:13"
$ws.Range("T13").Value = 6.17
$ws.Range("U13").Value = "23/09/2023 19:13"
$ws.Range("V13").Value = "https://www.betexplorer.com/football/algeria/ligue-1/biskra-us-souf/KYnDaKS7/"

# --- Append new rows 41..44 (copy formatting from row 40 first) ---
$ws.Range("A40:V40").Copy($ws.Range("A41:V41"))
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "algeria"
$ws.Range("C41").Value = "ligue-1"
$ws.Range("D41").Value = "2023-2024"
$ws.Range("E41").Value = 45247.64583333334
$ws.Range("F41").Value = "Constantine"
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = "Magra"
$ws.Range("I41").Value = 1
$ws.Range("J41").Value = 1.59
$ws.Range("K41").Value = "26/10/2023 04:42"
$ws.Range("L41").Value = 1.31
$ws.Range("M41").Value = "17/11/2023 15:24"
$ws.Range("N41").Value = 3.47
$ws.Range("O41").Value = "26/10/2023 04:42"
$ws.Range("P41").Value = 4.98
$ws.Range("Q41").Value = "17/11/2023 15:24"
$ws.Range("R41").Value = 5.42
$ws.Range("S41").Value = "26/10/2023 04:42"
$ws.Range("T41").Value = 10.79
$ws.Range("U41").Value = "17/11/2023 15:24"
$ws.Range("V41").Value = "https://www.betexplorer.com/football/algeria/ligue-1/constantine-magra/jgeqGfPg/"

$ws.Range("A40:V40").Copy($ws.Range("A42:V42"))
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "algeria"
$ws.Range("C42").Value = "ligue-1"
$ws.Range("D42").Value = "2023-2024"
$ws.Range("E42").Value = 45247.64583333334
$ws.Range("F42").Value = "El Bayadh"
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = "US Souf"
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 1.45
$ws.Range("K42").Value = "16/11/2023 03:42"
$ws.Range("L42").Value = 1.45
$ws.Range("M42").Value = "17/11/2023 15:18"
$ws.Range("N42").Value = 3.74
$ws.Range("O42").Value = "16/11/2023 03:42"
$ws.Range("P42").Value = 4.01
$ws.Range("Q42").Value = "17/11/2023 15:18"
$ws.Range("R42").Value = 7.38
$ws.Range("S42").Value = "16/11/2023 03:42"
$ws.Range("T42").Value = 8.609999999999999
$ws.Range("U42").Value = "17/11/2023 15:18"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/algeria/ligue-1/el-bayadh-us-souf/UqfmFEv0/"

$ws.Range("A40:V40").Copy($ws.Range("A43:V43"))
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "algeria"
$ws.Range("C43").Value = "ligue-1"
$ws.Range("D43").Value = "2023-2024"
$ws.Range("E43").Value = 45247.69791666666
$ws.Range("F43").Value = "Oran"
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = "Khenchela"
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2.65
$ws.Range("K43").Value = "16/11/2023 03:42"
$ws.Range("L43").Value = 2.27
$ws.Range("M43").Value = "17/11/2023 16:43"
$ws.Range("N43").Value = 2.95
$ws.Range("O43").Value = "16/11/2023 03:42"
$ws.Range("P43").Value = 2.83
$ws.Range("Q43").Value = "17/11/2023 16:43"
$ws.Range("R43").Value = 2.6
$ws.Range("S43").Value = "16/11/2023 03:42"
$ws.Range("T43").Value = 3.81
$ws.Range("U43").Value = "17/11/2023 16:43"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/algeria/ligue-1/oran-khenchela/dIttHzAm/"

$ws.Range("A40:V40").Copy($ws.Range("A44:V44"))
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "algeria"
$ws.Range("C44").Value = "ligue-1"
$ws.Range("D44").Value = "2023-2024"
$ws.Range("E44").Value = 45247.75
$ws.Range("F44").Value = "MC Alger"
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = "Saoura"
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 1.56
$ws.Range("K44").Value = "16/11/2023 06:12"
$ws.Range("L44").Value = 1.49
$ws.Range("M44").Value = "17/11/2023 17:56"
$ws.Range("N44").Value = 3.38
$ws.Range("O44").Value = "16/11/2023 06:12"
$ws.Range("P44").Value = 3.73
$ws.Range("Q44").Value = "17/11/2023 17:56"
$ws.Range("R44").Value = 5.99
$ws.Range("S44").Value = "16/11/2023 06:12"
$ws.Range("T44").Value = 8.619999999999999
$ws.Range("U44").Value = "17/11/2023 17:56"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/algeria/ligue-1/mc-alger-saoura/C8qhEYg6/"

